$wb = $excel.ActiveWorkbook

# This script re-applies a scheduled market-price / profit recalculation
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.
# Each row holds a single FFXIV leve: H=currentAveragePrice,
# I=currentAveragePriceNQ, J=currentAveragePriceHQ, K=LevePriceNQ,
# L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ.


$ws = $wb.Worksheets.Item("ALC")
# Row 20
$ws.Range("H20").Value = 4666.6665
$ws.Range("I20").Value = 4666.6665
$ws.Range("K20").Value = 4666.6665
$ws.Range("M20").Value = -4436.6665

# Row 31
$ws.Range("H31").Value = 99
$ws.Range("I31").Value = 99
$ws.Range("K31").Value = 297
$ws.Range("M31").Value = -67

# Row 35
$ws.Range("H35").Value = 4666.6665
$ws.Range("I35").Value = 4666.6665
$ws.Range("K35").Value = 4666.6665
$ws.Range("M35").Value = -4287.6665

# Row 38
$ws.Range("H38").Value = 1842.3572
$ws.Range("I38").Value = 224.875
$ws.Range("J38").Value = 3999
$ws.Range("K38").Value = 674.625
$ws.Range("L38").Value = 11997
$ws.Range("M38").Value = -302.625
$ws.Range("N38").Value = -12741

# Row 62
$ws.Range("H62").Value = 11299
$ws.Range("I62").Value = 7558.8
$ws.Range("K62").Value = 7558.8
$ws.Range("M62").Value = -6934.8

# Row 65
$ws.Range("H65").Value = 11299
$ws.Range("I65").Value = 7558.8
$ws.Range("K65").Value = 37794
$ws.Range("M65").Value = -34674

# Row 106
$ws.Range("H106").Value = 2438.95
$ws.Range("I106").Value = 2293.2778
$ws.Range("J106").Value = 3750
$ws.Range("K106").Value = 2293.2778
$ws.Range("L106").Value = 3750
$ws.Range("M106").Value = -1662.2778
$ws.Range("N106").Value = -5012

# Row 113
$ws.Range("H113").Value = 24162.5
$ws.Range("J113").Value = 6216.6665
$ws.Range("L113").Value = 6216.6665
$ws.Range("N113").Value = -12724.6665

# Row 116
$ws.Range("H116").Value = 4005.5715
$ws.Range("I116").Value = 3686.25
$ws.Range("J116").Value = 4431.3335
$ws.Range("K116").Value = 3686.25
$ws.Range("L116").Value = 4431.3335
$ws.Range("M116").Value = -244.25
$ws.Range("N116").Value = -11315.3335


$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 1000
$ws.Range("I5").Value = 1000
$ws.Range("K5").Value = 1000
$ws.Range("M5").Value = -888

# Row 45
$ws.Range("H45").Value = 1680.25
$ws.Range("I45").Value = 1578.6
$ws.Range("J45").Value = 1849.6666
$ws.Range("K45").Value = 1578.6
$ws.Range("L45").Value = 1849.6666
$ws.Range("M45").Value = -1201.6
$ws.Range("N45").Value = -2603.6666

# Row 46
$ws.Range("H46").Value = 7268.7144
$ws.Range("I46").Value = 6813.5
$ws.Range("K46").Value = 6813.5
$ws.Range("M46").Value = -6494.5

# Row 97
$ws.Range("H97").Value = 2320.0435
$ws.Range("I97").Value = 2576.25
$ws.Range("K97").Value = 2576.25
$ws.Range("M97").Value = -2080.25

# Row 118
$ws.Range("H118").Value = 250000
$ws.Range("J118").Value = 250000
$ws.Range("L118").Value = 250000
$ws.Range("N118").Value = -253314

# Row 122
$ws.Range("H122").Value = 1430658.6
$ws.Range("I122").Value = 2224158
$ws.Range("K122").Value = 6672474
$ws.Range("M122").Value = -6670024


$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 1000
$ws.Range("I4").Value = 1000
$ws.Range("K4").Value = 1000
$ws.Range("M4").Value = -885

# Row 42
$ws.Range("H42").Value = 399684
$ws.Range("J42").Value = 399684
$ws.Range("L42").Value = 399684
$ws.Range("N42").Value = -400340

# Row 63
$ws.Range("H63").Value = 24755
$ws.Range("J63").Value = 24755
$ws.Range("L63").Value = 24755
$ws.Range("N63").Value = -26127

# Row 66
$ws.Range("H66").Value = 24755
$ws.Range("J66").Value = 24755
$ws.Range("L66").Value = 74265
$ws.Range("N66").Value = -81129

# Row 86
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()

# Row 89
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()


$ws = $wb.Worksheets.Item("CRP")
# Row 86
$ws.Range("H86").Value = 9904.571
$ws.Range("I86").Value = 10842.777
$ws.Range("K86").Value = 10842.777
$ws.Range("M86").Value = -9719.777

# Row 89
$ws.Range("H89").Value = 9904.571
$ws.Range("I89").Value = 10842.777
$ws.Range("K89").Value = 54213.885
$ws.Range("M89").Value = -48597.885

# Row 99
$ws.Range("H99").Value = 3755
$ws.Range("I99").Value = 3887.647
$ws.Range("K99").Value = 3887.647
$ws.Range("M99").Value = -2389.647

# Row 126
$ws.Range("H126").Value = 3755
$ws.Range("I126").Value = 3887.647
$ws.Range("K126").Value = 11662.941
$ws.Range("M126").Value = -9192.940999999999


$ws = $wb.Worksheets.Item("CUL")
# Row 125
$ws.Range("H125").Value = 2000
$ws.Range("I125").Value = 2000
$ws.Range("K125").Value = 6000
$ws.Range("M125").Value = -1080

# Row 131
$ws.Range("H131").Value = 3194.8145
$ws.Range("J131").Value = 2627.117
$ws.Range("L131").Value = 7881.351000000001
$ws.Range("N131").Value = -17961.351

# Row 140
$ws.Range("H140").Value = 2292.9048


$ws = $wb.Worksheets.Item("GSM")
# Row 99
$ws.Range("H99").Value = 9735.5
$ws.Range("I99").Value = 9735.5
$ws.Range("K99").Value = 9735.5
$ws.Range("M99").Value = -7489.5

# Row 102
$ws.Range("H102").Value = 6057.0625
$ws.Range("I102").Value = 10614.143
$ws.Range("J102").Value = 2512.6667
$ws.Range("K102").Value = 10614.143
$ws.Range("L102").Value = 2512.6667
$ws.Range("M102").Value = -8992.143
$ws.Range("N102").Value = -5756.6667

# Row 122
$ws.Range("H122").Value = 2576.2727
$ws.Range("I122").Value = 2333.1
$ws.Range("J122").Value = 5008
$ws.Range("K122").Value = 6999.299999999999
$ws.Range("L122").Value = 15024
$ws.Range("M122").Value = -4549.299999999999
$ws.Range("N122").Value = -19924


$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 971.4545000000001
$ws.Range("I22").Value = 395.2857
$ws.Range("J22").Value = 1979.75
$ws.Range("K22").Value = 395.2857
$ws.Range("L22").Value = 1979.75
$ws.Range("M22").Value = -100.2857
$ws.Range("N22").Value = -2569.75

# Row 27
$ws.Range("H27").Value = 971.4545000000001
$ws.Range("I27").Value = 395.2857
$ws.Range("J27").Value = 1979.75
$ws.Range("K27").Value = 395.2857
$ws.Range("L27").Value = 1979.75
$ws.Range("M27").Value = -288.2857
$ws.Range("N27").Value = -2193.75

# Row 46
$ws.Range("H46").Value = 3554.889
$ws.Range("J46").Value = 3554.889
$ws.Range("L46").Value = 3554.889
$ws.Range("N46").Value = -3930.889

# Row 82
$ws.Range("H82").Value = 2075.0908
$ws.Range("I82").Value = 862.25
$ws.Range("J82").Value = 2768.1428
$ws.Range("K82").Value = 862.25
$ws.Range("L82").Value = 2768.1428
$ws.Range("M82").Value = -501.25
$ws.Range("N82").Value = -3490.1428

# Row 85
$ws.Range("H85").Value = 2075.0908
$ws.Range("I85").Value = 862.25
$ws.Range("J85").Value = 2768.1428
$ws.Range("K85").Value = 862.25
$ws.Range("L85").Value = 2768.1428
$ws.Range("M85").Value = 385.75
$ws.Range("N85").Value = -5264.1428

# Row 93
$ws.Range("H93").Value = 9209.105
$ws.Range("I93").Value = 8936.5625
$ws.Range("J93").Value = 10662.667
$ws.Range("K93").Value = 8936.5625
$ws.Range("L93").Value = 10662.667
$ws.Range("M93").Value = -7688.5625
$ws.Range("N93").Value = -13158.667

# Row 99
$ws.Range("H99").Value = 34556.3
$ws.Range("I99").Value = 26173.666
$ws.Range("K99").Value = 26173.666
$ws.Range("M99").Value = -23178.666

# Row 100
$ws.Range("H100").Value = 3178.45
$ws.Range("I100").Value = 2971.5334
$ws.Range("J100").Value = 3799.2
$ws.Range("K100").Value = 2971.5334
$ws.Range("L100").Value = 3799.2
$ws.Range("M100").Value = -2430.5334
$ws.Range("N100").Value = -4881.2

# Row 122
$ws.Range("H122").Value = 7738.9165
$ws.Range("I122").Value = 8299.200000000001
$ws.Range("K122").Value = 24897.6
$ws.Range("M122").Value = -22447.6

# Row 132
$ws.Range("H132").Value = 1072272.4
$ws.Range("I132").Value = 4525.9443
$ws.Range("K132").Value = 13577.8329
$ws.Range("M132").Value = -11047.8329


$ws = $wb.Worksheets.Item("WVR")
# Row 33
$ws.Range("H33").Value = 28999.334
$ws.Range("I33").Value = 18999
$ws.Range("K33").Value = 18999
$ws.Range("M33").Value = -18749

# Row 36
$ws.Range("H36").Value = 28999.334
$ws.Range("I36").Value = 18999
$ws.Range("K36").Value = 18999
$ws.Range("M36").Value = -18749

# Row 81
$ws.Range("H81").Value = 1484.5454
$ws.Range("I81").Value = 1247.5714
$ws.Range("J81").Value = 1899.25
$ws.Range("K81").Value = 2495.1428
$ws.Range("L81").Value = 3798.5
$ws.Range("M81").Value = -1434.1428
$ws.Range("N81").Value = -5920.5

# Row 84
$ws.Range("H84").Value = 1484.5454
$ws.Range("I84").Value = 1247.5714
$ws.Range("J84").Value = 1899.25
$ws.Range("K84").Value = 12475.714
$ws.Range("L84").Value = 18992.5
$ws.Range("M84").Value = -7171.714
$ws.Range("N84").Value = -29600.5

# Row 96
$ws.Range("H96").Value = 1215
$ws.Range("I96").Value = 672.5
$ws.Range("J96").Value = 2300
$ws.Range("K96").Value = 672.5
$ws.Range("L96").Value = 2300
$ws.Range("M96").Value = 700.5
$ws.Range("N96").Value = -5046
